$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: D2, E2
$ws.Range("D2").Value = "70.748.09"
$ws.Range("E2").Value = "  +1.55%  "

# Row 3: D3, E3
$ws.Range("D3").Value = "3.631.71"
$ws.Range("E3").Value = "  +3.60%  "

# Row 4: E4
$ws.Range("E4").Value = "  +0.08%  "

# Row 5: D5, E5
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "604.88"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.04%  "

# Row 6: D6, E6
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "198.99"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +2.12%  "

# Row 7: D7, E7
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "0.627"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.28%  "

# Row 9: D9, E9
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "0.221"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +9.92%  "

# Row 10: D10, E10
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "0.647"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -0.45%  "

# Row 11: D11, E11
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "53.87"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.77%  "

# Row 12: E12
$ws.Range("E12").Value = "  +1.87%  "

# Row 13: D13, E13
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "9.56"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +0.41%  "

# Row 14: D14, E14
$ws.Range("D14").Value = "4.209.09"
$ws.Range("E14").Value = "  +3.53%  "

# Row 15: D15, E15
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "677.34"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +14.02%  "

# Row 16: D16, E16
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "13.04"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +2.18%  "

# Row 17: B17, C17, D17, E17
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "70.881.32"
$ws.Range("E17").Value = "  +1.50%  "

# Row 18: B18, C18, D18, E18
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.636.22"
$ws.Range("E18").Value = "  +3.72%  "

# Row 19: D19, E19
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "19.08"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.26%  "

# Row 20: E20
$ws.Range("E20").Value = "  +0.40%  "

# Row 21: D21, E21
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +1.16%  "

# Row 22: D22, E22
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "18.50"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.84%  "

# Row 23: E23
$ws.Range("E23").Value = "  +1.51%  "

# Row 24: D24, E24
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "105.68"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +3.97%  "

# Row 25: D25, E25
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "4.62"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.36%  "

# Row 26: E26
$ws.Range("E26").Value = "  -5.05%  "

# Row 27: D27, E27
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "10.46"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -3.57%  "

# Row 28: D28, E28
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "9.82"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +2.98%  "

# Row 29: D29, E29
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "34.01"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +2.41%  "

# Row 30: D30, E30
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "4.66"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +8.64%  "

# Row 31: D31, E31
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "7.20"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +1.92%  "

# Row 32: D32
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "12.22"
$c.Style = "Normal"

# Row 33: D33, E33
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "0.115"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.41%  "

# Row 34: D34, E34
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "63.46"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.50%  "

# Row 35: D35, E35
$ws.Range("D35").Value = "3.973.04"
$ws.Range("E35").Value = "  +6.48%  "

# Row 36: D36, E36
$c = $ws.Cells.Item(36, 4)
$c.Value = "0.0PLACEHOLDER0867"
$c.Replace("PLACEHOLDER", [string][char]0x2083) | Out-Null
$ws.Range("E36").Value = "  +5.52%  "

# Row 37: D37, E37
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.06%  "

# Row 38: E38
$ws.Range("E38").Value = "  -1.61%  "

# Row 39: D39, E39
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "505.72"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +4.55%  "

# Row 40: D40
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "36.81"
$c.Style = "Normal"

# Row 41: E41
$ws.Range("E41").Value = "  -0.67%  "

# Row 42: D42, E42
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "3.54"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -3.35%  "

# Row 43: E43
$ws.Range("E43").Value = "  +2.50%  "

# Row 44: D44, E44
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "3.09"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +9.55%  "

# Row 45: D45, E45
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "0.0460"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +1.46%  "

# Row 46: D46, E46
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "3.51"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +6.90%  "

# Row 47: E47
$ws.Range("E47").Value = "  +0.74%  "

# Row 48: D48, E48
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "8.68"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +3.34%  "

# Row 49: E49
$ws.Range("E49").Value = "  -0.26%  "

# Row 50: D50, E50
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "0.000249"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +1.28%  "

# Row 51: D51, E51
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "2.97"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +5.21%  "
